# Rename the sheet from "Sheet2" to "Sheet1"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet1"

# Make the header row styling consistent across all four columns
# (A1 and C1 previously used a plain bold style; B1/D1 used bold + centered/top
# aligned. Copy B1's format onto A1 and C1 so the whole header row matches.)
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("C1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Replace the "Marketing_Amount" column (C2:C101) with the newly generated values
$c = New-Object 'object[,]' 100,1
$c[0,0] = 254.7026172983832
$c[1,0] = 92.75786041836116
$c[2,0] = 174.186899205287
$c[3,0] = 293.5446599600729
$c[4,0] = 373.1278995074983
$c[5,0] = 176.8861060061794
$c[6,0] = 231.2544208762795
$c[7,0] = 299.1821395851151
$c[8,0] = 39.08905741032211
$c[9,0] = 216.2799250969186
$c[10,0] = 252.4521785580439
$c[11,0] = 291.7136753481487
$c[12,0] = 294.5518862573497
$c[13,0] = 108.8337508246414
$c[14,0] = 112.1931616372713
$c[15,0] = 120.1837163687133
$c[16,0] = 409.2039536578803
$c[17,0] = 262.74208681171
$c[18,0] = 197.9033850825451
$c[19,0] = 24.04521303491376
$c[20,0] = -24.89949079170394
$c[21,0] = 200.680929772018
$c[22,0] = 146.7218099429753
$c[23,0] = 253.1417489796779
$c[24,0] = 274.7377311993804
$c[25,0] = -16.46828372993824
$c[26,0] = 289.5379258650723
$c[27,0] = 82.14080749870831
$c[28,0] = 312.8889607179229
$c[29,0] = 413.9679384950143
$c[30,0] = 131.4973712848458
$c[31,0] = 238.6581259801087
$c[32,0] = 140.8607126184943
$c[33,0] = 189.7101765888036
$c[34,0] = 343.3543925336924
$c[35,0] = 381.317448455199
$c[36,0] = 352.5145340363861
$c[37,0] = 171.8689924392206
$c[38,0] = 150.8836591296024
$c[39,0] = 204.6348624712332
$c[40,0] = 198.8223517466454
$c[41,0] = 76.74910314105124
$c[42,0] = 62.43649046874937
$c[43,0] = 433.5387697615895
$c[44,0] = 329.2673909124173
$c[45,0] = 231.5962849194407
$c[46,0] = 311.6102319975037
$c[47,0] = 275.1245177915955
$c[48,0] = 95.55510762210243
$c[49,0] = 122.1129859893016
$c[50,0] = 16.72667194031622
$c[51,0] = 80.0951248929631
$c[52,0] = 163.4597431215564
$c[53,0] = 161.4683907938794
$c[54,0] = 149.3408885830673
$c[55,0] = 181.1665935265209
$c[56,0] = 318.3258611191584
$c[57,0] = 189.8735948869891
$c[58,0] = 26.78389531595182
$c[59,0] = 236.8629417006431
$c[60,0] = 99.87697761120245
$c[61,0] = 115.5223419229729
$c[62,0] = 311.8426858977773
$c[63,0] = -9.06413011658384
$c[64,0] = 338.1213071126876
$c[65,0] = 45.91095318958691
$c[66,0] = 138.9900826516978
$c[67,0] = 265.3891127762887
$c[68,0] = 243.3850817808379
$c[69,0] = 331.0972697898069
$c[70,0] = 152.7045281088768
$c[71,0] = 75.44914553787054
$c[72,0] = 425.9700342271651
$c[73,0] = -0.9912910176826344
$c[74,0] = 339.6170820588774
$c[75,0] = 196.0094954529843
$c[76,0] = 234.7101425409059
$c[77,0] = 327.3075167617792
$c[78,0] = 64.67237339363137
$c[79,0] = 258.5582671114873
$c[80,0] = 247.7425079608322
$c[81,0] = 171.7913243477093
$c[82,0] = 261.033121986523
$c[83,0] = 283.1878156861388
$c[84,0] = 173.41260968978
$c[85,0] = 200.5444588015292
$c[86,0] = 384.4389785579825
$c[87,0] = 327.0037582093825
$c[88,0] = 58.21236892447288
$c[89,0] = 129.9725863465568
$c[90,0] = 314.341152651341
$c[91,0] = 189.3722535191214
$c[92,0] = 367.413748903843
$c[93,0] = 149.3319518241856
$c[94,0] = 207.5683198587201
$c[95,0] = 315.8286584095974
$c[96,0] = 234.525001036041
$c[97,0] = 162.0435246952918
$c[98,0] = 139.095604635181
$c[99,0] = 155.8494681722351
$ws.Range("C2:C101").Value = $c

# Resize columns B and C to fit their (now updated) content
$ws.Columns.Item(2).ColumnWidth = 5.5
$ws.Columns.Item(3).ColumnWidth = 16.666666666666668

# Update the active selection
$ws.Range("D3").Select() | Out-Null
